$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(6).Delete()
$ws.Rows.Item(35).Delete()
